$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.630.86"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "'1.755.61"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'324.28"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4505"
$ws.Range("E7").Value = "  +5.37%  "
$ws.Range("D8").Value = "'0.3555"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "'0.07508"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'41.60"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "'1.088"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "'5.994"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "'7.158"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "'1.754.26"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "'93.42"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "'0.06498"
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'17.05"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'5.752"
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("D23").Value = "'27.670.01"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").Value = "'2.109"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "'163.66"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "'1.957.26"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "'2.083"
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("D30").Value = "'125.33"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").Value = "'0.09185"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").Value = "'3.659"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").Value = "'5.498"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").Value = "'0.02285"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").Value = "'11.68"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("D37").Value = "'0.06036"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").Value = "'0.6294"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "'4.943"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "'1.182"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").Value = "'1.391"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "'7.765"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "'13.16"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").Value = "'3.708"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'0.5871"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "'123.19"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'1.943"
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("D49").Value = "'0.06898"
$ws.Range("D50").Value = "'1.132"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("E51").Value = "  -2.37%  "
